# Daily countdown update:
# For every data row (row 2 .. last used row) on the active sheet:
#   D = total days, E = remaining days, F = cycle start date (YYYYMMDD as a plain number)
#
# One more day has passed since the last update, so:
#   - if the remaining-days counter (E) is already down to 1, the cycle has
#     completed: restart it by advancing the start date (F) forward by the
#     full cycle length (D) and resetting the remaining days (E) back to D.
#   - otherwise simply decrement the remaining days (E) by 1.
#
# Rows whose start date (F) is not a clean 8-digit YYYYMMDD value are left
# untouched (the date is corrupt, so there is nothing sensible to roll
# forward).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {

    $d = $ws.Cells.Item($r, 4).Value()
    $e = $ws.Cells.Item($r, 5).Value()
    $f = $ws.Cells.Item($r, 6).Value()

    if ($d -eq $null -or $e -eq $null -or $f -eq $null) {
        continue
    }

    $fStr = [string]$f

    if ($fStr.Length -ne 8) {
        # Malformed start date - skip this row entirely.
        continue
    }

    $year = [int]$fStr.Substring(0, 4)
    $month = [int]$fStr.Substring(4, 2)
    $day = [int]$fStr.Substring(6, 2)

    if ($e -eq 1) {
        # Cycle finished - roll the start date forward by the full
        # cycle length and reset the remaining-days counter.
        $startDate = Get-Date -Year $year -Month $month -Day $day
        $newDate = $startDate.AddDays($d)

        $ws.Cells.Item($r, 6).Value = [int]$newDate.ToString("yyyyMMdd")
        $ws.Cells.Item($r, 5).Value = $d
    } else {
        # Still mid-cycle - one more day has elapsed.
        $ws.Cells.Item($r, 5).Value = $e - 1
    }
}
